$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 2.316462327490609, -12.90434384910116),
    @(3, 9.162319139560072, 11.00161740362515),
    @(4, 4.440675708811326, 0.02592873473603596),
    @(5, 6.166830030198267, 5.009734572843771),
    @(6, 2.204591531842581, -4.726949348040732),
    @(7, 1.09689030557385, -4.442195584720931),
    @(8, 1.600060471414833, -2.949819494134909),
    @(9, 1.651055586686678, 2.833734372666652),
    @(10, 1.99440460461342, 0.9898525518331924),
    @(11, 2.613530175870626, 3.313641510414356),
    @(12, 3.696331036365752, 6.233648892987009),
    @(13, -0.973865931199458, -4.327930935899992),
    @(14, 1.725130460355095, -1.194610791900008),
    @(15, -0.5115004854862049, 2.23509962177757),
    @(16, 0.9727820482463123, 2.866869504079239),
    @(17, 2.122313752051319, 0.8382457967197388),
    @(18, -0.1133034947815914, 0.8323378752418176),
    @(19, 3.161594928268019, 2.57979941834241)
)

foreach ($row in $values) {
    $r = $row[0]
    $cVal = $row[1]
    $eVal = $row[2]
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 5).Value = $eVal
}
